$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the new row of data (A2/B2) that was previously blank.
$ws.Range("A2").Value = "Rittmang"
$ws.Range("B2").Value = 8452047071

# Extend the formatted table body down to row 18 by copying the
# formatting of the last existing row (row 6) into the new rows,
# matching the row height used throughout the sheet.
$ws.Range("A6:B6").Copy()
$ws.Range("A7:B18").PasteSpecial(-4122)
$ws.Range("A7:B18").EntireRow.RowHeight = 15.75
$excel.CutCopyMode = $false
